$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Encode" section values (D4:D6)
$ws.Range("D4").Value = 0.234
$ws.Range("D5").Value = 2.312
$ws.Range("D6").Value = 48.906

# Update the "Query" section values (D18:D26)
$ws.Range("D18").Value = 0.001
$ws.Range("D19").Value = 0.001
$ws.Range("D20").Value = 0.001
$ws.Range("D21").Value = 0.001
$ws.Range("D22").Value = 0.001
$ws.Range("D23").Value = 0.001
$ws.Range("D24").Value = 0.016
$ws.Range("D25").Value = 0.031
$ws.Range("D26").Value = 0.091

# Update selected cell in the sheet view
$ws.Range("F7").Select()
